# Update the wcl_weights column (Q) from 6000 to 2000 for all data rows
# that currently hold the value 6000 (rows 2-18 and 21-73 in this sheet;
# rows 19-20 hold different computed weight values and are left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 17).End(-4162).Row  # xlUp = -4162, column 17 = Q

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 17)  # column Q
    if ($cell.Value2 -eq 6000) {
        $cell.Value2 = 2000
    }
}
